# [BI-1059] Updating files for name length
# Rename the trait name cells and update method-description text on the
# "Template" worksheet of empty_then_2_rows.xlsx

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Template")

# Row 3 (leaf trait): shorten the *Trait name and update *Method name wording
$ws.Range("A3").Value = "PM_Leaf"
# Row 4 (stalk trait): shorten the *Trait name
$ws.Range("A4").Value = "PM_Stalk"

# Update the method-name text for both rows
$ws.Range("J3").Value = "Powdery Mildew severity, leaf"
$ws.Range("J4").Value = "Powdery Mildew severity, stalk"

# Refresh the active selection on the sheet (was S5, now J5; the scrolled
# topLeftCell is reset back to the default as well)
$ws.Activate()
$ws.Range("J5").Select()
